$wb = $excel.ActiveWorkbook

# --- Sheet: "Fatalities by Age Group" ---
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsAge.Range("B4").Value = 77
$wsAge.Range("B5").Value = 650
$wsAge.Range("B6").Value = 2081
$wsAge.Range("B7").Value = 5083
$wsAge.Range("B8").Value = 9670
$wsAge.Range("B9").Value = 7372
$wsAge.Range("B10").Value = 8675
$wsAge.Range("B11").Value = 9212
$wsAge.Range("B12").Value = 8700
$wsAge.Range("B13").Value = 20527
$wsAge.Range("B15").Value = 72082

# --- Sheet: "Fatalities by Gender" ---
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsGender.Range("B2").Value = 30176
$wsGender.Range("B3").Value = 41905
[void]$wsGender.Select()
[void]$wsGender.Range("B2:B4").Select()

# --- Sheet: "Fatalities by Race-Ethnicity" ---
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsRace.Range("B2").Value = 1333
$wsRace.Range("B3").Value = 7630
$wsRace.Range("B4").Value = 31482
$wsRace.Range("B5").Value = 428
$wsRace.Range("B6").Value = 31164
[void]$wsRace.Select()
[void]$wsRace.Range("B8").Select()
